$d = $word.ActiveDocument
$paras = $d.Paragraphs
$offset = 0

# paragraph 11: 'Colorado state funding volatility for higher educa'
$p = $paras.Item(11 + $offset)
$p.Range.Text = "Trump administration (2025–2029) reducing federal HE funding; 120 TRIO programs terminated"

# paragraph 12: 'Federal financial aid policy changes (Pell Grant, '
$p = $paras.Item(12 + $offset)
$p.Range.Text = "DEI programs under HIGH scrutiny — executive order targeting DEI in accreditation (Apr 2025)"

# paragraph 13: 'Native American tuition waiver mandate (federal ob'
$p = $paras.Item(13 + $offset)
$p.Range.Text = "Tribal education funding VOLATILE: 109% increase Sept 2025, but FY2026 proposes 24% cuts"

# paragraph 14: 'State performance-based funding models'
$p = $paras.Item(14 + $offset)
$p.Range.Text = "Colorado FY 2025–26: `$38.4M increase (far less than `$95M requested); 3.5% tuition cap"

# paragraph 15: 'Political pressure on DEI programs in public insti'
$p = $paras.Item(15 + $offset)
$p.Range.Text = "Native American Tuition Waiver at risk of misclassification as DEI (waiver is statutory, not DEI)`r"
$offset = $offset + 1
$p2 = $paras.Item(15 + $offset)
$p2.Range.Text = "HLC providing flexibility on diversity standards, but federal pressure on accreditors continues"

# paragraph 17: 'Leverage federal tribal education funding'
$p = $paras.Item(17 + $offset)
$p.Range.Text = "Reframe Indigenous programs through statutory obligations (CRS 23-52-105) and cultural preservation (legally safe)`r"
$offset = $offset + 1
$p2 = $paras.Item(17 + $offset)
$p2.Range.Text = "Use 'first-generation support' and 'inclusive excellence' framing (avoids identity-based language)"

# paragraph 19: 'Economic (Impact: High, Trend: Mixed)'
$p = $paras.Item(19 + $offset)
$p.Range.Text = "Economic (Impact: High, Trend: Negative)"

# paragraph 20: 'Declining state appropriations per student'
$p = $paras.Item(20 + $offset)
$p.Range.Text = "Colorado shifts costs to students via tuition rather than state appropriations"

# paragraph 21: 'Rising tuition sensitivity among families'
$p = $paras.Item(21 + $offset)
$p.Range.Text = "Rising tuition sensitivity; students increasingly price-conscious and comparison-shopping"

# paragraph 22: 'Durango cost of living affecting faculty recruitme'
$p = $paras.Item(22 + $offset)
$p.Range.Text = "Durango housing crisis — major hidden barrier for student attendance AND faculty recruitment"

# paragraph 23: 'Native American tuition waiver revenue impact (~37'
$p = $paras.Item(23 + $offset)
$p.Range.Text = "Native American tuition waiver revenue impact (~37% of students at zero tuition)"

# paragraph 24: 'Economic diversification in Four Corners region'
$p = $paras.Item(24 + $offset)
$p.Range.Text = "Regional economy tourism-dependent (seasonal, variable); limited large employers"

# paragraph 25: 'Student debt burden concerns nationally'
$p = $paras.Item(25 + $offset)
$p.Range.Text = "Skills-based hiring growing — degrees less of an automatic hiring requirement"

# paragraph 27: 'Grow graduate programs for additional revenue'
$p = $paras.Item(27 + $offset)
$p.Range.Text = "Healthcare/nursing programs (strong regional employer demand)"

# paragraph 28: 'Expand dual enrollment pipeline'
$p = $paras.Item(28 + $offset)
$p.Range.Text = "Expand dual enrollment pipeline (Pueblo CC, San Juan College feeders)"

# paragraph 29: 'Develop workforce-aligned certificates'
$p = $paras.Item(29 + $offset)
$p.Range.Text = "Develop workforce-aligned certificates and micro-credentials`r"
$offset = $offset + 1
$p2 = $paras.Item(29 + $offset)
$p2.Range.Text = "Position as affordable rural alternative to cost-climbing urban institutions"

# paragraph 30: 'Social (Impact: High, Trend: Mixed)'
$p = $paras.Item(30 + $offset)
$p.Range.Text = "Social (Impact: Medium-High, Trend: Mixed)"

# paragraph 31: 'Declining college-going rates nationally'
$p = $paras.Item(31 + $offset)
$p.Range.Text = "Declining college-going rates nationally and in Colorado"

# paragraph 32: 'Changing student expectations (career-focused outc'
$p = $paras.Item(32 + $offset)
$p.Range.Text = "Career outcome expectations dominant ('What job will I get?')"

# paragraph 33: 'Growing demand for flexible/hybrid learning'
$p = $paras.Item(33 + $offset)
$p.Range.Text = "Indigenous education opportunity IS REAL (166 tribes, 37% waiver, underserved nationally)"

# paragraph 34: 'FLC unique mission serving Native American student'
$p = $paras.Item(34 + $offset)
$p.Range.Text = "First-generation students (43%) need targeted support systems"

# paragraph 35: 'First-generation students (43%) need additional su'
$p = $paras.Item(35 + $offset)
$p.Range.Text = "Growing skepticism about ROI of higher education; trade/vocational paths gaining acceptance"

# paragraph 36: 'Mental health and wellness demands increasing'
$p = $paras.Item(36 + $offset)
$p.Range.Text = "Strong outdoor/recreation culture aligns with FLC place-based brand"

# paragraph 38: 'Outdoor recreation lifestyle as recruitment differ'
$p = $paras.Item(38 + $offset)
$p.Range.Text = "Indigenous education leadership — reframe through statutory obligations (CRS 23-52-105), not DEI"

# paragraph 39: 'Indigenous education leadership positioning'
$p = $paras.Item(39 + $offset)
$p.Range.Text = "First-generation student success programs (safe framing, encompasses many Indigenous students)"

# paragraph 40: 'Experiential learning emphasis'
$p = $paras.Item(40 + $offset)
$p.Range.Text = "Place-based brand leveraging Durango outdoor lifestyle as recruitment differentiator`r"
$offset = $offset + 1
$p2 = $paras.Item(40 + $offset)
$p2.Range.Text = "Career outcome emphasis across all programs"

# paragraph 41: 'Technological (Impact: Medium, Trend: Opportunity)'
$p = $paras.Item(41 + $offset)
$p.Range.Text = "Technological (Impact: High, Trend: Rapidly Changing)"

# paragraph 42: 'AI disruption in curriculum and pedagogy'
$p = $paras.Item(42 + $offset)
$p.Range.Text = "AI disruption transforming pedagogy, assessment, and student expectations"

# paragraph 43: 'Need for technology infrastructure upgrades'
$p = $paras.Item(43 + $offset)
$p.Range.Text = "Online graduate market SATURATED — ASU, SNHU, Western Governors dominate (`$50M+ marketing)"

# paragraph 44: 'Online/hybrid program delivery expectations'
$p = $paras.Item(44 + $offset)
$p.Range.Text = "FLC has NO online brand nationally; ~25 online courses (~10% of offerings)"

# paragraph 45: 'Data analytics for student success and retention'
$p = $paras.Item(45 + $offset)
$p.Range.Text = "Passive video lectures becoming obsolete; AI-enabled adaptive learning replacing them"

# paragraph 46: 'AI Institute at FLC as emerging strength'
$p = $paras.Item(46 + $offset)
$p.Range.Text = "AI Institute at FLC as emerging institutional strength`r"
$offset = $offset + 1
$p2 = $paras.Item(46 + $offset)
$p2.Range.Text = "Online program development requires 1–2+ years governance + substantial investment"

# paragraph 48: 'AI Institute partnerships and growth'
$p = $paras.Item(48 + $offset)
$p.Range.Text = "AI Institute partnerships and curriculum integration"

# paragraph 49: 'Technology-enhanced experiential learning'
$p = $paras.Item(49 + $offset)
$p.Range.Text = "AI-enabled advising, early alerts, and retention prediction tools"

# paragraph 50: 'Online graduate program expansion'
$p = $paras.Item(50 + $offset)
$p.Range.Text = "AI literacy across all disciplines as differentiator"

# paragraph 51: 'Legal (Impact: Medium, Trend: Stable)'
$p = $paras.Item(51 + $offset)
$p.Range.Text = "Legal (Impact: High, Trend: Deteriorating)"

# paragraph 52: 'Accreditation compliance requirements (HLC)'
$p = $paras.Item(52 + $offset)
$p.Range.Text = "Title VI scrutiny — 50+ universities under investigation for race-conscious programs"

# paragraph 53: 'Title IX and student safety regulations'
$p = $paras.Item(53 + $offset)
$p.Range.Text = "Native American Tuition Waiver has DISTINCT legal basis (CRS 23-52-105, since 1911)"

# paragraph 54: 'Federal reporting mandates (IPEDS)'
$p = $paras.Item(54 + $offset)
$p.Range.Text = "HLC accreditation: federal pressure on DEI standards, but HLC offers flexibility"

# paragraph 55: 'Employment law for faculty/staff'
$p = $paras.Item(55 + $offset)
$p.Range.Text = "Trump administration revising Title IX regulations (definitions, due process in flux)"

# paragraph 56: 'Tribal sovereignty considerations in partnerships'
$p = $paras.Item(56 + $offset)
$p.Range.Text = "FERPA compliance critical for AI tools processing student data`r"
$offset = $offset + 1
$p2 = $paras.Item(56 + $offset)
$p2.Range.Text = "Programs framed as 'equity-focused' are primary federal targets"

# paragraph 58: 'Streamlined accreditation through proactive compli'
$p = $paras.Item(58 + $offset)
$p.Range.Text = "NATW defensible under Title VI (statutory basis per CRS 23-52-105, not voluntary DEI)"

# paragraph 59: 'Tribal education partnership agreements'
$p = $paras.Item(59 + $offset)
$p.Range.Text = "Government-to-government tribal partnerships (sovereignty framing, not race-based)`r"
$offset = $offset + 1
$p2 = $paras.Item(59 + $offset)
$p2.Range.Text = "HLC flexibility allows alternative methods to meet diversity-related standards"

# paragraph 60: 'Environmental (Impact: Medium, Trend: Opportunity)'
$p = $paras.Item(60 + $offset)
$p.Range.Text = "Environmental (Impact: Medium, Trend: Negative)"

# paragraph 61: 'Climate change impacts on Durango/mountain region'
$p = $paras.Item(61 + $offset)
$p.Range.Text = "Southwest Colorado wildfire risk increasing — smoke impacts air quality and outdoor activities"

# paragraph 62: 'Campus sustainability expectations from students'
$p = $paras.Item(62 + $offset)
$p.Range.Text = "Colorado River basin under long-term drought stress; water rights contentious"

# paragraph 63: 'Environmental science as program strength'
$p = $paras.Item(63 + $offset)
$p.Range.Text = "Snowpack variability affects regional economy (ski, rafting, outdoor recreation)"

# paragraph 64: 'Outdoor recreation economy dependency on climate'
$p = $paras.Item(64 + $offset)
$p.Range.Text = "Outdoor recreation brand is FLC strength but CLIMATE-VULNERABLE"

# paragraph 65: 'Wildfire risk to campus and community'
$p = $paras.Item(65 + $offset)
$p.Range.Text = "Sustainability compliance is baseline, not differentiator"

# paragraph 67: 'Position as leader in sustainability education'
$p = $paras.Item(67 + $offset)
$p.Range.Text = "Proactive sustainability initiatives to build brand beyond compliance"

# paragraph 68: 'Climate resilience research opportunities'
$p = $paras.Item(68 + $offset)
$p.Range.Text = "Emergency preparedness planning as operational strength"

# paragraph 69: 'Green campus initiatives for recruitment'
$p = $paras.Item(69 + $offset)
$p.Range.Text = "Environmental science/conservation programs align with regional needs"

# paragraph 72: 'Technological, Legal, and Environmental factors ar'
$p = $paras.Item(72 + $offset)
$p.Range.Text = "Technological factors (4/5 impact) reflect a rapidly changing landscape: the online market is saturated and FLC has no online brand, though AI is transforming delivery. Legal factors (4/5) are deteriorating with NATW legal basis (CRS 23-52-105) requiring proactive documentation. Environmental factors present both risk (wildfire/drought) and opportunity (outdoor brand)."

# paragraph 74: 'The PESTLE analysis reveals that FLC operates in a'
$p = $paras.Item(74 + $offset)
$p.Range.Text = "The PESTLE analysis reveals that FLC operates in an environment of heightened political risk and constrained resources. Federal DEI policy disruptions threaten TRIO programs and could misclassify the statutory NATW mission. The Durango housing crisis constrains faculty/staff recruitment. Key strategic imperatives include: (1) reframing Indigenous education through statutory obligations (CRS 23-52-105) and state law (not DEI), (2) strengthening the dual enrollment pipeline as a hedge against declining first-year enrollment, (3) improving retention as the most cost-effective enrollment strategy, and (4) investing in AI capabilities while recognizing that large-scale online expansion faces a saturated market where FLC has no brand."

# paragraph 76: 'Prioritize revenue diversification through graduat'
$p = $paras.Item(76 + $offset)
$p.Range.Text = "Proactively document NATW statutory basis (CRS 23-52-105) to protect against DEI misclassification"

# paragraph 77: 'Strengthen advocacy for state funding while reduci'
$p = $paras.Item(77 + $offset)
$p.Range.Text = "Prioritize retention improvement as the most cost-effective enrollment strategy (Compass, early-alert)"

# paragraph 78: 'Position Indigenous education mission as a nationa'
$p = $paras.Item(78 + $offset)
$p.Range.Text = "Grow dual enrollment pipeline and transfer pathways as near-term enrollment stabilizers"

# paragraph 79: 'Invest in AI Institute and sustainability programs'
$p = $paras.Item(79 + $offset)
$p.Range.Text = "Frame Indigenous education through statutory obligations (CRS 23-52-105) and cultural preservation, not DEI language"

# paragraph 80: 'Develop workforce-aligned certificates and micro-c'
$p = $paras.Item(80 + $offset)
$p.Range.Text = "Invest in AI Institute and experiential learning as place-based institutional differentiators"

# paragraph 81: 'Build data-driven retention programs targeting equ'
$p = $paras.Item(81 + $offset)
$p.Range.Text = "Qualify online expansion: pursue only Indigenous niche (NATW moat), not generic online degrees`r"
$offset = $offset + 1
$p2 = $paras.Item(81 + $offset)
$p2.Range.Text = "Address Durango housing crisis impact on faculty/staff recruitment through institutional partnerships"
